# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 28;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 60;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 71;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 73;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 75;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 80;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 85;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 87;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 90;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 101; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 113; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 129; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 133; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 136; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 137; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 143; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 148; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 160; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 161; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 166; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 167; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 212; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 216; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 220; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 224; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 234; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 235; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 297; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 305; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 307; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 318; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 326; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 334; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
